$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 6 data (Trial 4)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "April 21 2024"
$ws.Range("C6").Value = "0.21, 0.085, 0.115, 0.526"
$ws.Range("D6").Value = 0.215
$ws.Range("E6").Value = "6pi"
$ws.Range("F6").Value = 418.9
$ws.Range("G6").Value = "not saved"
$ws.Range("H6").Value = "same as above"
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = "CPU"
$ws.Range("K6").Value = "Pass"

# New header for column L, then the new data point
$ws.Range("L1").Value = "Time [s]"
$ws.Range("L6").Value = 1035.9739999999999

# Update selection to match the recorded final state
$ws.Range("H13").Select()
